$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 4 for columns D, J, K, L, M, O, P
# (these are the columns whose values differ between the two rows / change in the diff)

# Row 2 -> new values (previously held by row 4)
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("O2").Value = "Limache"
$ws.Range("P2").Value = 183

# Row 4 -> new values (previously held by row 2)
$ws.Range("D4").Value = 44273
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 233
